# Weekly update: insert a new week of price data (Primera / Segunda) for
# "Betarraga" at the top of the data block (rows 1048-1049), pushing the
# existing records down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right above the current row 1048, shifting all
# the existing data (old rows 1048:1139) down to rows 1050:1141.
$ws.Range("A1048:A1049").EntireRow.Insert()

# --- New row 1048 (Primera) ---
$ws.Range("A1048").Value = 3
$ws.Range("B1048").Value = "Femacal de La Calera"
$ws.Range("C1048").Value = "Coquimbo"
$ws.Range("D1048").Value = 45106
$ws.Range("E1048").Value = 5
$ws.Range("F1048").Value = 100114014
$ws.Range("G1048").Value = "Betarraga"
$ws.Range("H1048").Value = "Sin especificar"
$ws.Range("I1048").Value = "Primera"
$ws.Range("J1048").Value = 5400
$ws.Range("K1048").Value = 500
$ws.Range("L1048").Value = 550
$ws.Range("M1048").Value = 525
$ws.Range("N1048").Value = "`$/paquete 4 unidades"
$ws.Range("O1048").Value = "Provincia de Quillota"
$ws.Range("P1048").Value = 131
$ws.Range("Q1048").Value = 4
$ws.Range("R1048").Value = "Hortaliza"

# --- New row 1049 (Segunda) ---
$ws.Range("A1049").Value = 3
$ws.Range("B1049").Value = "Femacal de La Calera"
$ws.Range("C1049").Value = "Coquimbo"
$ws.Range("D1049").Value = 45106
$ws.Range("E1049").Value = 5
$ws.Range("F1049").Value = 100114014
$ws.Range("G1049").Value = "Betarraga"
$ws.Range("H1049").Value = "Sin especificar"
$ws.Range("I1049").Value = "Segunda"
$ws.Range("J1049").Value = 3700
$ws.Range("K1049").Value = 400
$ws.Range("L1049").Value = 420
$ws.Range("M1049").Value = 410
$ws.Range("N1049").Value = "`$/paquete 4 unidades"
$ws.Range("O1049").Value = "Provincia de Quillota"
$ws.Range("P1049").Value = 102
$ws.Range("Q1049").Value = 4
$ws.Range("R1049").Value = "Hortaliza"
